# annotate script, reorder chunks
# Renames raw column headers (column B term labels) to human-readable
# labels, and refreshes the computed p-values (column C) across every
# histology-result worksheet.

$wb = $excel.ActiveWorkbook

# --- 1. Human-readable term labels -----------------------------------
$termRenames = @{
    "age"                       = "Age at Diagnosis"
    "reported_gender"           = "Reported Gender"
    "race"                      = "Reported Race"
    "ethnicity"                 = "Reported Ethnicity"
    "mol_sub_group"             = "Molecular Subgroup"
    "extent_of_tumor_resection" = "Extent of Tumor Resection"
}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 2)
        $label = $cell.Value()
        if ($termRenames.ContainsKey($label)) {
            $cell.Value = $termRenames[$label]
        }
    }
}

# --- 2. Updated p-values (column C) per worksheet ---------------------
$valueChanges = @(
    @("Low-grade glioma", "C2", 0.217969285650775),
    @("Low-grade glioma", "C3", 0.606339366063394),
    @("Low-grade glioma", "C6", 0.0331966803319668),
    @("Low-grade glioma", "C7", 0.228377162283772),
    @("Low-grade glioma", "C8", 0.787121287871213),
    @("Low-grade glioma", "C9", 0.0612938706129387),
    @("Low-grade glioma", "C10", 0.0438762953173008),
    @("Ependymoma", "C2", 0.278583212588568),
    @("Ependymoma", "C3", 0.405059494050595),
    @("Ependymoma", "C6", 0.262573742625737),
    @("Ependymoma", "C7", 0.278672132786721),
    @("Ependymoma", "C8", 0.966403359664034),
    @("Ependymoma", "C11", 0.743616658580255),
    @("DIPG or DMG", "C3", 0.505149485051495),
    @("DIPG or DMG", "C6", 0.0158984101589841),
    @("DIPG or DMG", "C7", 0.618238176182382),
    @("DIPG or DMG", "C9", 0.124845238053223),
    @("DIPG or DMG", "C10", 0.259262165005025),
    @("ATRT", "C3", 0.368663133686631),
    @("ATRT", "C5", 0.0007999200079992),
    @("ATRT", "C6", 0.498550144985501),
    @("ATRT", "C7", 0.284671532846715),
    @("ATRT", "C8", 0.171182881711829),
    @("ATRT", "C9", 0.212678732126787),
    @("ATRT", "C11", 0.541788422608163),
    @("Other high-grade glioma", "C3", 0.7996200379962),
    @("Other high-grade glioma", "C6", 0.837116288371163),
    @("Other high-grade glioma", "C7", 0.0172982701729827),
    @("Other high-grade glioma", "C8", 0.668633136686331),
    @("Other high-grade glioma", "C9", 0.975602439756024),
    @("Other high-grade glioma", "C10", 0.415013092936939),
    @("Meningioma", "C2", 0.0928659156743261),
    @("Meningioma", "C3", 0.942105789421058),
    @("Meningioma", "C5", 0.0003999600039996),
    @("Meningioma", "C6", 0.68953104689531),
    @("Meningioma", "C7", 0.311968803119688),
    @("Meningioma", "C8", 0.201379862013799),
    @("Neurofibroma plexiform", "C4", 0.0002999700029997),
    @("Neurofibroma plexiform", "C5", 0.149185081491851),
    @("Neurofibroma plexiform", "C6", 0.0546945305469453),
    @("Neurofibroma plexiform", "C7", 0.294770522947705),
    @("Neurofibroma plexiform", "C8", 0.398860113988601),
    @("Oligodendroglioma", "C2", 0.0575209855960223),
    @("Oligodendroglioma", "C3", 0.300569943005699),
    @("Oligodendroglioma", "C5", 0.291870812918708),
    @("Oligodendroglioma", "C7", 0.221877812218778),
    @("Non-neoplastic tumor", "C2", 0.563235079302386),
    @("Non-neoplastic tumor", "C3", 0.905509449055094),
    @("Non-neoplastic tumor", "C5", 0.0001999800019998),
    @("Non-neoplastic tumor", "C6", 0.250874912508749),
    @("Non-neoplastic tumor", "C7", 0.854214578542146),
    @("Mixed neuronal-glial tumor", "C3", 0.201379862013799),
    @("Mixed neuronal-glial tumor", "C6", 0.521047895210479),
    @("Mixed neuronal-glial tumor", "C7", 0.784521547845215),
    @("Mixed neuronal-glial tumor", "C8", 0.710928907109289),
    @("Mixed neuronal-glial tumor", "C9", 0.261273872612739),
    @("Medulloblastoma", "C2", 0.890643284295059),
    @("Medulloblastoma", "C3", 0.783321667833217),
    @("Medulloblastoma", "C6", 0.477352264773523),
    @("Medulloblastoma", "C7", 0.297670232976702),
    @("Medulloblastoma", "C8", 0.836516348365164),
    @("Schwannoma", "C3", 0.913008699130087),
    @("Schwannoma", "C5", 0.520847915208479),
    @("Schwannoma", "C6", 0.925207479252075),
    @("Schwannoma", "C7", 0.724727527247275),
    @("Schwannoma", "C8", 0.944905509449055),
    @("Schwannoma", "C10", 0.0419697378487613),
    @("Mesenchymal tumor", "C3", 0.126887311268873),
    @("Mesenchymal tumor", "C6", 0.841915808419158),
    @("Mesenchymal tumor", "C7", 0.484551544845515),
    @("Mesenchymal tumor", "C8", 0.166383361663834),
    @("Mesenchymal tumor", "C10", 0.588256408874304),
    @("Germ cell tumor", "C3", 0.628537146285371),
    @("Germ cell tumor", "C5", 0.0076992300769923),
    @("Germ cell tumor", "C6", 0.369163083691631),
    @("Germ cell tumor", "C7", 0.277872212778722),
    @("Germ cell tumor", "C9", 0.344188269692038),
    @("Craniopharyngioma", "C3", 0.733826617338266),
    @("Craniopharyngioma", "C5", 0.0038996100389961),
    @("Craniopharyngioma", "C6", 0.165083491650835),
    @("Craniopharyngioma", "C7", 0.58004199580042),
    @("Other tumor", "C2", 0.744389813120937),
    @("Other tumor", "C3", 0.0970902909709029),
    @("Other tumor", "C5", 0.0005999400059994),
    @("Other tumor", "C6", 0.71002899710029),
    @("Other tumor", "C9", 0.715220319009891)
)

foreach ($change in $valueChanges) {
    $sheetName = $change[0]
    $cellRef = $change[1]
    $newValue = $change[2]
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range($cellRef).Value = $newValue
}
